$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7034338116645813
$ws.Range("B1").Value = 3.904030799865723
$ws.Range("C1").Value = 5.734167575836182
$ws.Range("D1").Value = 1.23979640007019
$ws.Range("E1").Value = 0.7141823768615723
